# UC007 - Listar Autorizações de Pagamento Pendentes
# Bump version 1.0 -> 1.2.5 and apply the wording/content fixes from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version bump
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text fix (appears once per test case block: TC1..TC5)
$preconditionText = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value  = $preconditionText
$ws.Range("B17").Value = $preconditionText
$ws.Range("B25").Value = $preconditionText
$ws.Range("B33").Value = $preconditionText
$ws.Range("B41").Value = $preconditionText

# 3) TC1 step wording fixes (accentuation + trailing periods)
$ws.Range("D10").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B11").Value = "Chefe Seleciona uma diária apta para pagamento."
$ws.Range("D11").Value = "SYSTEM Destaca a diária selecionada."

# 4) TC2/TC3/TC4 steps were reordered: the "filter" pair moves into the TC2 slot,
#    the "register payment authorization" pair moves into the TC3 slot, and the
#    "assign/unassign" pair moves into the TC4 slot (with a small wording tweak).
$ws.Range("B20").Value = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$ws.Range("D20").Value = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

$ws.Range("B28").Value = "Chefe Clica para realizar a autorização de pagamento."
$ws.Range("D28").Value = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento."

$ws.Range("B36").Value = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D36").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# 5) TC5 step wording fix (trailing period)
$ws.Range("D44").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
